$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.11811066666667
$ws.Range("H2").Value = 138.354332
$ws.Range("I2").Value = 0.95896098489411
$ws.Range("J2").Value = 0.9589609848941099
$ws.Range("M2").Value = 0.02179466666666667
$ws.Range("N2").Value = 0.065384
$ws.Range("O2").Value = 0.06643270670809397
$ws.Range("P2").Value = 0.06643270670809397
$ws.Range("Q2").Value = 1.005128849276445
$ws.Range("R2").Value = 9.046159643488
$ws.Range("S2").Value = 0.06370637385397533
$ws.Range("T2").Value = 0.06370637385397533

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.11811066666667
$ws.Range("H3").Value = 138.354332
$ws.Range("I3").Value = 0.95896098489411
$ws.Range("J3").Value = 0.9589609848941099
$ws.Range("M3").Value = 0.2375286666666667
$ws.Range("N3").Value = 0.7125860000000001
$ws.Range("O3").Value = 0.7240153056144294
$ws.Range("P3").Value = 0.7240153056144294
$ws.Range("Q3").Value = 10.95437333583911
$ws.Range("R3").Value = 98.589360022552
$ws.Range("S3").Value = 0.6943024305504233
$ws.Range("T3").Value = 0.6943024305504232

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.11811066666667
$ws.Range("H4").Value = 138.354332
$ws.Range("I4").Value = 0.95896098489411
$ws.Range("J4").Value = 0.9589609848941099
$ws.Range("O4").Value = 0.2095519876774766
$ws.Range("P4").Value = 0.2095519876774767
$ws.Range("Q4").Value = 3.170527872112
$ws.Range("R4").Value = 28.534750849008
$ws.Range("S4").Value = 0.2009521804897114
$ws.Range("T4").Value = 0.2009521804897114

$ws.Range("I5").Value = 0.002799731840346333
$ws.Range("J5").Value = 0.002799731840346333
$ws.Range("M5").Value = 0.02179466666666667
$ws.Range("N5").Value = 0.065384
$ws.Range("O5").Value = 0.06643270670809397
$ws.Range("P5").Value = 0.06643270670809397
$ws.Range("Q5").Value = 0.002934521098666666
$ws.Range("R5").Value = 0.026410689888
$ws.Range("S5").Value = 0.0001859937642110401
$ws.Range("T5").Value = 0.0001859937642110401

$ws.Range("I6").Value = 0.002799731840346333
$ws.Range("J6").Value = 0.002799731840346333
$ws.Range("M6").Value = 0.2375286666666667
$ws.Range("N6").Value = 0.7125860000000001
$ws.Range("O6").Value = 0.7240153056144294
$ws.Range("P6").Value = 0.7240153056144294
$ws.Range("S6").Value = 0.002027048704026799
$ws.Range("T6").Value = 0.002027048704026799

$ws.Range("I7").Value = 0.002799731840346333
$ws.Range("J7").Value = 0.002799731840346333
$ws.Range("O7").Value = 0.2095519876774766
$ws.Range("P7").Value = 0.2095519876774767
$ws.Range("Q7").Value = 0.009256505711999999
$ws.Range("R7").Value = 0.08330855140799999
$ws.Range("S7").Value = 0.0005866893721084937
$ws.Range("T7").Value = 0.0005866893721084938

$ws.Range("I8").Value = 0.0382392832655437
$ws.Range("J8").Value = 0.0382392832655437
$ws.Range("M8").Value = 0.02179466666666667
$ws.Range("N8").Value = 0.065384
$ws.Range("O8").Value = 0.06643270670809397
$ws.Range("P8").Value = 0.06643270670809397
$ws.Range("Q8").Value = 0.04008026123200001
$ws.Range("R8").Value = 0.360722351088
$ws.Range("S8").Value = 0.00254033908990759
$ws.Range("T8").Value = 0.00254033908990759

$ws.Range("I9").Value = 0.0382392832655437
$ws.Range("J9").Value = 0.0382392832655437
$ws.Range("M9").Value = 0.2375286666666667
$ws.Range("N9").Value = 0.7125860000000001
$ws.Range("O9").Value = 0.7240153056144294
$ws.Range("P9").Value = 0.7240153056144294
$ws.Range("Q9").Value = 0.4368137928280001
$ws.Range("R9").Value = 3.931324135452001
$ws.Range("S9").Value = 0.02768582635997936
$ws.Range("T9").Value = 0.02768582635997936

$ws.Range("I10").Value = 0.0382392832655437
$ws.Range("J10").Value = 0.0382392832655437
$ws.Range("O10").Value = 0.2095519876774766
$ws.Range("P10").Value = 0.2095519876774767
$ws.Range("S10").Value = 0.008013117815656751
$ws.Range("T10").Value = 0.008013117815656752
